$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'242.31"
$ws.Range("D3").Value = "'23.38"
$ws.Range("D4").Value = "'5.645"
$ws.Range("D5").Value = "'0.05814"
$ws.Range("D6").Value = "'3.413"
$ws.Range("D7").Value = "'6.469"
$ws.Range("D8").Value = "'1.317"
$ws.Range("D9").Value = "'0.7981"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01243"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1457"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07618"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03217"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.02967"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09238"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001663"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "MCDex"
$ws.Range("C17").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D17").Value = "'3.332"
$ws.Range("E17").Value = "16MCDexMCB"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04750"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("D19").Value = "'0.006201"
$ws.Range("D20").Value = "'0.005458"
$ws.Range("D21").Value = "'0.001069"
$ws.Range("D23").Value = "'3.694"
$ws.Range("D25").Value = "'0.3325"
$ws.Range("D27").Value = "'0.001001"
$ws.Range("D40").Value = "'0.04283"
$ws.Range("D41").Value = "'0.007107"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003603"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1053"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.009530"
$ws.Range("E45").Value = "44ACDXExchangeACXTWorstin24h"
$ws.Range("D46").Value = "'0.00005434"
$ws.Range("D49").Value = "'0.1047"
$ws.Range("E49").Value = "48BOLOBOLO"
$ws.Range("D50").Value = "'0.00002102"